# Apply crypto price/volume updates to Sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.317.34"

$ws.Cells.Item(3, 4).Value = "1.875.11"
$ws.Cells.Item(3, 5).Value = "  +0.82%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.7124"
$ws.Cells.Item(5, 5).Value = "  -0.17%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "241.93"
$ws.Cells.Item(6, 5).Value = "  +0.64%  "

$ws.Cells.Item(7, 5).Value = "  +0.02%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3102"
$ws.Cells.Item(8, 5).Value = "  +0.84%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07733"
$ws.Cells.Item(9, 5).Value = "  -0.38%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "24.86"
$ws.Cells.Item(10, 5).Value = "  -1.05%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.08517"
$ws.Cells.Item(11, 5).Value = "  +3.11%  "

$ws.Cells.Item(12, 4).Value = "1.888.00"
$ws.Cells.Item(12, 5).Value = "  +1.03%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.217"
$ws.Cells.Item(13, 5).Value = "  -0.27%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.7098"
$ws.Cells.Item(14, 5).Value = "  -0.91%  "

$ws.Cells.Item(15, 5).Value = "  +1.22%  "

$ws.Cells.Item(16, 4).Value = "29.312.96"
$ws.Cells.Item(16, 5).Value = "  +0.38%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.000008200"
$ws.Cells.Item(17, 5).Value = "  +5.18%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "6.001"
$ws.Cells.Item(18, 5).Value = "  +2.37%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "241.59"
$ws.Cells.Item(19, 5).Value = "  -1.24%  "

$ws.Cells.Item(20, 4).Value = "2.134.12"
$ws.Cells.Item(20, 5).Value = "  +1.20%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.23"
$ws.Cells.Item(21, 5).Value = "  +0.70%  "

$ws.Cells.Item(22, 5).Value = "  -0.03%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "7.799"
$ws.Cells.Item(23, 5).Value = "  -2.29%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.000"
$ws.Cells.Item(24, 5).Value = "  -0.03%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.1603"
$ws.Cells.Item(25, 5).Value = "  +0.73%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "163.02"
$ws.Cells.Item(26, 5).Value = "  +0.36%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "9.037"
$ws.Cells.Item(27, 5).Value = "  +1.45%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "18.47"
$ws.Cells.Item(28, 5).Value = "  +0.71%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.396"
$ws.Cells.Item(30, 5).Value = "  -0.57%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.317"
$ws.Cells.Item(31, 5).Value = "  +1.94%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.286"
$ws.Cells.Item(32, 5).Value = "  -2.25%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.05261"
$ws.Cells.Item(33, 5).Value = "  +1.46%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.931"
$ws.Cells.Item(34, 5).Value = "  +1.00%  "

$ws.Cells.Item(35, 5).Value = "  +0.24%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.7440"
$ws.Cells.Item(36, 5).Value = "  +2.30%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.686"
$ws.Cells.Item(37, 5).Value = "  +0.48%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01866"
$ws.Cells.Item(38, 5).Value = "  +0.61%  "

$ws.Cells.Item(40, 4).Value = "1.182.82"
$ws.Cells.Item(40, 5).Value = "  +1.48%  "

$ws.Cells.Item(41, 5).Value = "  +3.67%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.8885"
$ws.Cells.Item(42, 5).Value = "  -2.02%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "72.88"
$ws.Cells.Item(43, 5).Value = "  +0.65%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "106.26"
$ws.Cells.Item(44, 5).Value = "  +4.52%  "

$ws.Cells.Item(45, 5).Value = "  +0.04%  "

$ws.Cells.Item(46, 4).Value = "2.029.71"
$ws.Cells.Item(46, 5).Value = "  +1.37%  "

$ws.Cells.Item(48, 5).Value = "  -0.22%  "

$ws.Cells.Item(49, 5).Value = "  +1.06%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "9.382"
$ws.Cells.Item(50, 5).Value = "  +0.78%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.4313"
$ws.Cells.Item(51, 5).Value = "  +1.18%  "
